$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.948619544506073
$ws.Range("B1").Value = 2.977652072906494
$ws.Range("C1").Value = 4.24641752243042
$ws.Range("D1").Value = 2.084782600402832
$ws.Range("E1").Value = 1.241719245910645
